$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1718.8235
$ws.Range("J17").Value = 1718.8235
$ws.Range("L17").Value = 5156.470499999999
$ws.Range("N17").Value = -5492.470499999999

$ws.Range("H31").Value = 333333340
$ws.Range("I31").Value = 333333340
$ws.Range("K31").Value = 1000000020
$ws.Range("M31").Value = -999999790

$ws.Range("H55").Value = 610
$ws.Range("I55").Value = 387.5
$ws.Range("K55").Value = 387.5
$ws.Range("M55").Value = -173.5

$ws.Range("H80").Value = 3821.0715
$ws.Range("J80").Value = 3343.6667
$ws.Range("L80").Value = 10031.0001
$ws.Range("N80").Value = -12027.0001

$ws.Range("H83").Value = 3821.0715
$ws.Range("J83").Value = 3343.6667
$ws.Range("L83").Value = 30093.0003
$ws.Range("N83").Value = -40077.0003

$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

$ws.Range("H98").Value = 910.6
$ws.Range("I98").Value = 688.3570999999999
$ws.Range("J98").Value = 1429.1666
$ws.Range("K98").Value = 688.3570999999999
$ws.Range("L98").Value = 1429.1666
$ws.Range("M98").Value = 809.6429000000001
$ws.Range("N98").Value = -4425.1666

$ws.Range("H106").Value = 26557.143
$ws.Range("I106").Value = 26557.143
$ws.Range("K106").Value = 26557.143
$ws.Range("M106").Value = -25926.143

$ws.Range("H116").Value = 8900
$ws.Range("J116").Value = 8900
$ws.Range("L116").Value = 8900
$ws.Range("N116").Value = -15784

$ws.Range("H122").Value = 910.6
$ws.Range("I122").Value = 688.3570999999999
$ws.Range("J122").Value = 1429.1666
$ws.Range("K122").Value = 2065.0713
$ws.Range("L122").Value = 4287.4998
$ws.Range("M122").Value = 384.9287000000004
$ws.Range("N122").Value = -9187.4998

$ws.Range("H125").Value = 107147330
$ws.Range("J125").Value = 83338360
$ws.Range("L125").Value = 750045240
$ws.Range("N125").Value = -750050160

$ws.Range("H127").Value = 1967.75
$ws.Range("I127").Value = 2660
$ws.Range("K127").Value = 7980
$ws.Range("M127").Value = -3020

$ws.Range("H132").Value = 1235.5454
$ws.Range("I132").Value = 1229.2
$ws.Range("K132").Value = 3687.6
$ws.Range("M132").Value = -1157.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1830.4286
$ws.Range("I97").Value = 1779.25
$ws.Range("K97").Value = 1779.25
$ws.Range("M97").Value = -1283.25

$ws.Range("H122").Value = 3158.1667
$ws.Range("I122").Value = 2869.8
$ws.Range("K122").Value = 8609.400000000001
$ws.Range("M122").Value = -6159.400000000001

$ws.Range("H132").Value = 2926
$ws.Range("I132").Value = 3010.2727
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 9030.8181
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -6500.8181
$ws.Range("N132").Value = -11057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 895.8889
$ws.Range("I20").Value = 830
$ws.Range("J20").Value = 978.25
$ws.Range("K20").Value = 830
$ws.Range("L20").Value = 978.25
$ws.Range("M20").Value = -583
$ws.Range("N20").Value = -1472.25

$ws.Range("H86").Value = 1947.8
$ws.Range("I86").Value = 1913.3334
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 1913.3334
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = -790.3334
$ws.Range("N86").Value = -4245.5

$ws.Range("H89").Value = 1947.8
$ws.Range("I89").Value = 1913.3334
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 9566.666999999999
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = -3950.666999999999
$ws.Range("N89").Value = -21229.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 552.2857
$ws.Range("I22").Value = 522.125
$ws.Range("K22").Value = 522.125
$ws.Range("M22").Value = -172.125

$ws.Range("H31").Value = 1742.4

$ws.Range("H34").Value = 1742.4

$ws.Range("H86").Value = 10600.25
$ws.Range("I86").Value = 12332.667
$ws.Range("J86").Value = 9560.799999999999
$ws.Range("K86").Value = 12332.667
$ws.Range("L86").Value = 9560.799999999999
$ws.Range("M86").Value = -11209.667
$ws.Range("N86").Value = -11806.8

$ws.Range("H89").Value = 10600.25
$ws.Range("I89").Value = 12332.667
$ws.Range("J89").Value = 9560.799999999999
$ws.Range("K89").Value = 61663.335
$ws.Range("L89").Value = 47804
$ws.Range("M89").Value = -56047.335
$ws.Range("N89").Value = -59036

$ws.Range("H107").Value = 722.1429000000001
$ws.Range("I107").Value = 720.1539
$ws.Range("K107").Value = 720.1539
$ws.Range("M107").Value = 1199.8461

$ws.Range("H120").Value = 39985
$ws.Range("J120").Value = 39985
$ws.Range("L120").Value = 39985
$ws.Range("N120").Value = -47243

$ws.Range("H132").Value = 4450
$ws.Range("I132").Value = 5500
$ws.Range("J132").Value = 2350
$ws.Range("K132").Value = 16500
$ws.Range("L132").Value = 7050
$ws.Range("M132").Value = -13970
$ws.Range("N132").Value = -12110

$ws.Range("H134").Value = 2222
$ws.Range("I134").Value = 2222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6666
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("M134").Value = -4131

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 245724.33
$ws.Range("I2").Value = 220076.2
$ws.Range("J2").Value = 277784.5
$ws.Range("K2").Value = 1320457.2
$ws.Range("L2").Value = 1666707
$ws.Range("M2").Value = -1320344.2
$ws.Range("N2").Value = -1666933

$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -1027

$ws.Range("H38").Value = 40.444443
$ws.Range("I38").Value = 49
$ws.Range("K38").Value = 147
$ws.Range("M38").Value = 200

$ws.Range("H97").Value = 444.66666
$ws.Range("I97").Value = 444.6
$ws.Range("J97").Value = 444.75
$ws.Range("K97").Value = 1333.8
$ws.Range("L97").Value = 1334.25
$ws.Range("M97").Value = -837.8000000000002
$ws.Range("N97").Value = -2326.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1752216.9
$ws.Range("I3").Value = 1833400
$ws.Range("J3").Value = 1671033.6
$ws.Range("K3").Value = 1833400
$ws.Range("L3").Value = 1671033.6
$ws.Range("M3").Value = -1833284
$ws.Range("N3").Value = -1671265.6

$ws.Range("H40").Value = 12500
$ws.Range("I40").Value = 12000
$ws.Range("K40").Value = 12000
$ws.Range("M40").Value = -11849

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H113").Value = 865.3333
$ws.Range("I113").Value = 865.3333
$ws.Range("K113").Value = 865.3333
$ws.Range("M113").Value = 1304.6667

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7466.1665
$ws.Range("I7").Value = 3352
$ws.Range("K7").Value = 3352
$ws.Range("M7").Value = -3240

$ws.Range("H22").Value = 800.8
$ws.Range("I22").Value = 891.3333
$ws.Range("J22").Value = 665
$ws.Range("K22").Value = 891.3333
$ws.Range("L22").Value = 665
$ws.Range("M22").Value = -596.3333
$ws.Range("N22").Value = -1255

$ws.Range("H27").Value = 800.8
$ws.Range("I27").Value = 891.3333
$ws.Range("J27").Value = 665
$ws.Range("K27").Value = 891.3333
$ws.Range("L27").Value = 665
$ws.Range("M27").Value = -784.3333
$ws.Range("N27").Value = -879

$ws.Range("H40").Value = 3188.4443
$ws.Range("J40").Value = 5500
$ws.Range("L40").Value = 5500
$ws.Range("N40").Value = -5772

$ws.Range("H61").Value = 3677
$ws.Range("I61").Value = 3677
$ws.Range("K61").Value = 3677
$ws.Range("M61").Value = -3475

$ws.Range("H113").Value = 3677
$ws.Range("I113").Value = 3677
$ws.Range("K113").Value = 3677
$ws.Range("M113").Value = -1507

$ws.Range("H126").Value = 7466.1665
$ws.Range("I126").Value = 3352
$ws.Range("K126").Value = 10056
$ws.Range("M126").Value = -7586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2161.0833
$ws.Range("I136").Value = 2004.0526
$ws.Range("K136").Value = 6012.1578
$ws.Range("M136").Value = -3462.1578
